$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the leading "Dato" column (B); everything shifts one column left
$ws.Range("B1").EntireColumn.Delete()

# Remove now-unused trailing columns (previously BB1:BI1, now shifted to BA1:BH1)
$ws.Range("BA1:BH1").EntireColumn.Delete()

# Re-write header labels (renames, insertions, corrections) for columns B..AZ
$ws.Range('B1').Value = 'Navn'
$ws.Range('C1').Value = 'CPR nr.'
$ws.Range('D1').Value = 'HDM_kur_nr'
$ws.Range('E1').Value = 'Alder'
$ws.Range('F1').Value = 'Højde'
$ws.Range('G1').Value = 'Vægt'
$ws.Range('H1').Value = 'Overflade'
$ws.Range('I1').Value = 'Treatment_time_t0'
$ws.Range('J1').Value = 'Total_væskemængde'
$ws.Range('K1').Value = 'Mindst_2000_væskemængde'
$ws.Range('L1').Value = 'Plasma_kreatin_før_start'
$ws.Range('M1').Value = 'Durise_600ml_6timer'
$ws.Range('N1').Value = 'Furosemid'
$ws.Range('O1').Value = 'Forhydrering_6000ml_4timer'
$ws.Range('P1').Value = 'Dosis_natriumcarbonat_ved_lav_pH'
$ws.Range('Q1').Value = 'Sygeplejerske_navn_forhydering'
$ws.Range('R1').Value = 'Sygeplejerske_tid_forhydering'
$ws.Range('S1').Value = 'one_to_ten_MTX_dose'
$ws.Range('T1').Value = 'Sygeplejerske_navn_one_to_ten_MTX_dose'
$ws.Range('U1').Value = 'Sygeplejerske_tid_one_to_ten_MTX_dose'
$ws.Range('V1').Value = 'Kontinuerlig_infusion_start'
$ws.Range('W1').Value = 'nine_to_ten_MTX_dose'
$ws.Range('X1').Value = 'total_volume_MTX_and_hydration_liquid'
$ws.Range('Y1').Value = 'Sygeplejerske_navn_nine_to_ten_MTX_dose'
$ws.Range('Z1').Value = 'Sygeplejerske_tid_nine_to_ten_MTX_dose'
$ws.Range('AA1').Value = 'Hydreringsvæske_reduceret'
$ws.Range('AB1').Value = 'Se_MTX_t23'
$ws.Range('AC1').Value = 'P_kreatin_t23'
$ws.Range('AD1').Value = 'Se_MTX_t36'
$ws.Range('AE1').Value = 'P_kreatin_t36'
$ws.Range('AF1').Value = 'Hydreing_ved_høj_P_kreatin_t36'
$ws.Range('AG1').Value = 'Durise_ved_høj_P_kreatin_t36'
$ws.Range('AH1').Value = 'Hydreing_ved_normal_P_kreatin_t36'
$ws.Range('AI1').Value = 'Durise_ved_normal_P_kreatin_t36'
$ws.Range('AJ1').Value = 'Se_MTX_t42'
$ws.Range('AK1').Value = 'Første_dosis_calciumfolinat_t42'
$ws.Range('AL1').Value = 'Sygeplejerske_navn_Første_dosis_calciumfolinat_t42'
$ws.Range('AM1').Value = 'Sygeplejerske_tid_Første_dosis_calciumfolinat_t42'
$ws.Range('AN1').Value = 'Se_MTX_t48'
$ws.Range('AO1').Value = 'Anden_dosis_calciumfolinat_t48'
$ws.Range('AP1').Value = 'Sygeplejerske_navn_Anden_dosis_calciumfolinat_t48'
$ws.Range('AQ1').Value = 'Sygeplejerske_tid_Anden_dosis_calciumfolinat_t48'
$ws.Range('AR1').Value = 'Se_MTX_t54'
$ws.Range('AS1').Value = 'Tredje_dosis_calciumfolinat_t54'
$ws.Range('AT1').Value = 'Sygeplejerske_navn_Tredje_dosis_calciumfolinat_t54'
$ws.Range('AU1').Value = 'Sygeplejerske_tid_Tredje_dosis_calciumfolinat_t54'
$ws.Range('AV1').Value = 'Se_MTX_t66'
$ws.Range('AW1').Value = 'P_kreatin_t66'
$ws.Range('AX1').Value = 'Fjerde_dosis_calciumfolinat_t66'
$ws.Range('AY1').Value = 'Sygeplejerske_navn_Fjerde_dosis_calciumfolinat_t66'
$ws.Range('AZ1').Value = 'Sygeplejerske_tid_Fjerde_dosis_calciumfolinat_t66'
